$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1779.5217
$ws.Range("I40").Value = 1635.8
$ws.Range("J40").Value = 1819.4445
$ws.Range("K40").Value = 1635.8
$ws.Range("L40").Value = 1819.4445
$ws.Range("M40").Value = -1460.8
$ws.Range("N40").Value = -2169.4445

$ws.Range("H107").Value = 713.26666
$ws.Range("I107").Value = 699.9167
$ws.Range("J107").Value = 766.6667
$ws.Range("K107").Value = 699.9167
$ws.Range("L107").Value = 766.6667
$ws.Range("M107").Value = 1220.0833
$ws.Range("N107").Value = -4606.6667

$ws.Range("H125").Value = 1087.6875
$ws.Range("I125").Value = 633
$ws.Range("J125").Value = 1239.25
$ws.Range("K125").Value = 5697
$ws.Range("L125").Value = 11153.25
$ws.Range("M125").Value = -3237
$ws.Range("N125").Value = -16073.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1258.8572
$ws.Range("I2").Value = 977
$ws.Range("J2").Value = 2950
$ws.Range("K2").Value = 977
$ws.Range("L2").Value = 2950
$ws.Range("M2").Value = -864
$ws.Range("N2").Value = -3176

$ws.Range("H32").Value = 17844.102
$ws.Range("I32").Value = 19302.115
$ws.Range("K32").Value = 19302.115
$ws.Range("M32").Value = -19015.115

$ws.Range("H97").Value = 423.75
$ws.Range("I97").Value = 323.8
$ws.Range("J97").Value = 590.3333
$ws.Range("K97").Value = 323.8
$ws.Range("L97").Value = 590.3333
$ws.Range("M97").Value = 172.2
$ws.Range("N97").Value = -1582.3333

$ws.Range("H116").Value = 1258.8572
$ws.Range("I116").Value = 977
$ws.Range("J116").Value = 2950
$ws.Range("K116").Value = 977
$ws.Range("L116").Value = 2950
$ws.Range("M116").Value = 1317
$ws.Range("N116").Value = -7538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1258.8572
$ws.Range("I3").Value = 977
$ws.Range("J3").Value = 2950
$ws.Range("K3").Value = 977
$ws.Range("L3").Value = 2950
$ws.Range("M3").Value = -863
$ws.Range("N3").Value = -3178

$ws.Range("H80").Value = 712.46155
$ws.Range("I80").Value = 949.75
$ws.Range("J80").Value = 607
$ws.Range("K80").Value = 949.75
$ws.Range("L80").Value = 607
$ws.Range("M80").Value = 48.25
$ws.Range("N80").Value = -2603

$ws.Range("H83").Value = 712.46155
$ws.Range("I83").Value = 949.75
$ws.Range("J83").Value = 607
$ws.Range("K83").Value = 4748.75
$ws.Range("L83").Value = 3035
$ws.Range("M83").Value = 243.25
$ws.Range("N83").Value = -13019

$ws.Range("H99").Value = 700.55
$ws.Range("I99").Value = 507.69232
$ws.Range("J99").Value = 1058.7142
$ws.Range("K99").Value = 507.69232
$ws.Range("L99").Value = 1058.7142
$ws.Range("M99").Value = 990.30768
$ws.Range("N99").Value = -4054.7142

$ws.Range("H105").Value = 2912.88
$ws.Range("I105").Value = 2926.5557
$ws.Range("K105").Value = 2926.5557
$ws.Range("M105").Value = -1179.5557

$ws.Range("H134").Value = 39785.184
$ws.Range("I134").Value = 55152.633
$ws.Range("J134").Value = 3287.5
$ws.Range("K134").Value = 165457.899
$ws.Range("L134").Value = 9862.5
$ws.Range("M134").Value = -162922.899
$ws.Range("N134").Value = -14932.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11113718
$ws.Range("I31").Value = 2531.6924
$ws.Range("J31").Value = 40002804
$ws.Range("K31").Value = 2531.6924
$ws.Range("L31").Value = 40002804
$ws.Range("M31").Value = -2236.6924
$ws.Range("N31").Value = -40003394

$ws.Range("H34").Value = 11113718
$ws.Range("I34").Value = 2531.6924
$ws.Range("J34").Value = 40002804
$ws.Range("K34").Value = 2531.6924
$ws.Range("L34").Value = 40002804
$ws.Range("M34").Value = -2329.6924
$ws.Range("N34").Value = -40003208

$ws.Range("H58").Value = 853.8261
$ws.Range("I58").Value = 869.8889
$ws.Range("J58").Value = 796
$ws.Range("K58").Value = 869.8889
$ws.Range("L58").Value = 796
$ws.Range("M58").Value = -666.8889
$ws.Range("N58").Value = -1202

$ws.Range("H136").Value = 853.8261
$ws.Range("I136").Value = 869.8889
$ws.Range("J136").Value = 796
$ws.Range("K136").Value = 2609.6667
$ws.Range("L136").Value = 2388
$ws.Range("M136").Value = -59.66670000000022
$ws.Range("N136").Value = -7488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 55000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 55000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 165000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -166996

$ws.Range("H78").Value = 55000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 55000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 495000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -504984

$ws.Range("H80").Value = 2750
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -564

$ws.Range("H83").Value = 2750
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 180

$ws.Range("H131").Value = 1329345.8
$ws.Range("I131").Value = 10327.272
$ws.Range("J131").Value = 3402089.2
$ws.Range("K131").Value = 30981.816
$ws.Range("L131").Value = 10206267.6
$ws.Range("M131").Value = -25941.816
$ws.Range("N131").Value = -10216347.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 767.6539
$ws.Range("I97").Value = 752.2222
$ws.Range("K97").Value = 752.2222
$ws.Range("M97").Value = -256.2222

$ws.Range("H113").Value = 22729570
$ws.Range("I113").Value = 25002376
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 25002376
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -25000206
$ws.Range("N113").Value = -5840

$ws.Range("H122").Value = 3429.85
$ws.Range("I122").Value = 4819.4
$ws.Range("J122").Value = 2966.6667
$ws.Range("K122").Value = 14458.2
$ws.Range("L122").Value = 8900.000100000001
$ws.Range("M122").Value = -12008.2
$ws.Range("N122").Value = -13800.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2340
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -990

$ws.Range("H27").Value = 2340
$ws.Range("J27").Value = 400
$ws.Range("L27").Value = 400
$ws.Range("N27").Value = -614

$ws.Range("H40").Value = 5500
$ws.Range("I40").Value = 5500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5364
$ws.Range("N40").ClearContents()

$ws.Range("H46").Value = 1488.4546
$ws.Range("I46").Value = 650
$ws.Range("J46").Value = 1620.8422
$ws.Range("K46").Value = 650
$ws.Range("L46").Value = 1620.8422
$ws.Range("M46").Value = -462
$ws.Range("N46").Value = -1996.8422

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2055.1765
$ws.Range("J81").Value = 2088.4285
$ws.Range("L81").Value = 4176.857
$ws.Range("N81").Value = -6298.857

$ws.Range("H84").Value = 2055.1765
$ws.Range("J84").Value = 2088.4285
$ws.Range("L84").Value = 20884.285
$ws.Range("N84").Value = -31492.285

$ws.Range("H107").Value = 216.66667
$ws.Range("I107").Value = 216.66667
$ws.Range("J107").Value = 216.66667
$ws.Range("K107").Value = 650.00001
$ws.Range("L107").Value = 650.00001
$ws.Range("M107").Value = 1269.99999
$ws.Range("N107").Value = -4490.00001
